$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two changed indicator header strings (Gini / extreme poverty date ranges) ---
$ws.Range('J2').Value = 'Coefficient de Gini (estimations de la Banque Mondiale, données les plus récentes 2013-22)'
$ws.Range('K2').Value = '% de la population en extrême pauvreté ($2.15 par jour, mesure la plus récente 2013-22)(1)'

# --- Update data cells (columns D:K) across rows 3-98 per refreshed source data ---
$ws.Range('D3').Value = [double]"224655.61203334801"
$ws.Range('E3').Value = [double]"74790.564022787497"
$ws.Range('F3').Value = [double]"1.09200170875581"
$ws.Range('G3').Value = [double]"6511.0446015009002"
$ws.Range('H3').Value = [double]"25.765270173944401"
$ws.Range('I3').Value = [double]"76.459000000000003"
$ws.Range('D4').Value = [double]"42567.414795943499"
$ws.Range('E4').Value = [double]"18767.689364003199"
$ws.Range('F4').Value = [double]"11.8402084705577"
$ws.Range('G4').Value = [double]"16410.113303886901"
$ws.Range('H4').Value = [double]"6.6790352504637598"
$ws.Range('I4').Value = [double]"62.686"
$ws.Range('D5').Value = [double]"11652.1206325117"
$ws.Range('E5').Value = [double]"4743.6239432546699"
$ws.Range('F5').Value = [double]"7.8833065177182897"
$ws.Range('G5').Value = [double]"10230.314320788701"
$ws.Range('H5').Value = [double]"3.7187546606876101"
$ws.Range('I5').Value = [double]"50.375"
$ws.Range('D6').Value = [double]"5994.5296549919703"
$ws.Range('E6').Value = [double]"2549.7325902869002"
$ws.Range('F6').Value = [double]"2.1016989564100101"
$ws.Range('G6').Value = [double]"2872.4333799196902"
$ws.Range('H6').Value = [double]"6.0477459195062302"
$ws.Range('I6').Value = [double]"64.144000000000005"
$ws.Range('D7').Value = [double]"33341.730621460898"
$ws.Range('E7').Value = [double]"12463.357531965999"
$ws.Range('F7').Value = [double]"4.5667438343159503"
$ws.Range('G7').Value = [double]"1552.64059002466"
$ws.Range('I7').Value = [double]"68.27"
$ws.Range('D8').Value = [double]"43421.836052005703"
$ws.Range('E8').Value = [double]"15776.7509356145"
$ws.Range('G8').Value = [double]"1353.6720574747501"
$ws.Range('H8').Value = [double]"5.6884869912667204"
$ws.Range('I8').Value = [double]"78.088999999999999"
$ws.Range('D9').Value = [double]"25511.716660585102"
$ws.Range('E9').Value = [double]"12307.431273381801"
$ws.Range('F9').Value = [double]"2.6560508199969801"
$ws.Range('G9').Value = [double]"10003.708165701801"
$ws.Range('H9').Value = [double]"3.6202833870034001"
$ws.Range('I9').Value = [double]"58.167000000000002"
$ws.Range('D10').Value = [double]"872413.66486146301"
$ws.Range('E10').Value = [double]"418907.28902418399"
$ws.Range('G10').Value = [double]"14505.6559343808"
$ws.Range('H10').Value = [double]"4.5556051808843199"
$ws.Range('I10').Value = [double]"54.7"
$ws.Range('D11').Value = [double]"69550.892485829594"
$ws.Range('E11').Value = [double]"22147.649568093901"
$ws.Range('F11').Value = [double]"4.5987340444093396"
$ws.Range('G11').Value = [double]"3571.6348806793799"
$ws.Range('H11').Value = [double]"22.0207676245778"
$ws.Range('I11').Value = [double]"60.87"
$ws.Range('D12').Value = [double]"36196.036178664697"
$ws.Range('E12').Value = [double]"35967.375627026799"
$ws.Range('F12').Value = [double]"8.4619725508283494"
$ws.Range('G12').Value = [double]"2336.3884208104701"
$ws.Range('H12').Value = [double]"98.546224939953802"
$ws.Range('I12').Value = [double]"65.165000000000006"
$ws.Range('D13').Value = [double]"1365305.5539768001"
$ws.Range('E13').Value = [double]"618421.46388059901"
$ws.Range('F13').Value = [double]"4.4585256814049199"
$ws.Range('G13').Value = [double]"7128.9265582077396"
$ws.Range('H13').Value = [double]"4.5556051808843199"
$ws.Range('I13').Value = [double]"63.892499999999998"
$ws.Range('D14').Value = [double]"9895.0807500333503"
$ws.Range('E14').Value = [double]"3350.8099274818301"
$ws.Range('F14').Value = [double]"3.11930234606745"
$ws.Range('G14').Value = [double]"809.00319170509999"
$ws.Range('H14').Value = [double]"8.3134892842137393"
$ws.Range('I14').Value = [double]"78.486999999999995"
$ws.Range('D15').Value = [double]"111381.27962556299"
$ws.Range('E15').Value = [double]"45391.496408315303"
$ws.Range('F15').Value = [double]"3.6487176659821801"
$ws.Range('G15').Value = [double]"4091.2506508188399"
$ws.Range('H15').Value = [double]"2.2619561486631099"
$ws.Range('I15').Value = [double]"71.177999999999997"
$ws.Range('D16').Value = [double]"5045.4236760543699"
$ws.Range('E16').Value = [double]"2584.6193586337699"
$ws.Range('F16').Value = [double]"0.98292343122549997"
$ws.Range('G16').Value = [double]"1025.49261708422"
$ws.Range('H16').Value = [double]"4.2643644595907997"
$ws.Range('I16').Value = [double]"69.822000000000003"
$ws.Range('D17').Value = [double]"27162.187522451499"
$ws.Range('E17').Value = [double]"11800.1917091056"
$ws.Range('F17').Value = [double]"-1.1022521776400001"
$ws.Range('G17').Value = [double]"1605.8047604168801"
$ws.Range('H17').Value = [double]"-0.77284437035910003"
$ws.Range('I17').Value = [double]"59.314999999999998"
$ws.Range('D18').Value = [double]"21866.3640426235"
$ws.Range('E18').Value = [double]"12130.855358092"
$ws.Range('F18').Value = [double]"1.50003178327904"
$ws.Range('G18').Value = [double]"4556.6173218047697"
$ws.Range('H18').Value = [double]"1.9690254939855401"
$ws.Range('I18').Value = [double]"67.180999999999997"
$ws.Range('D19').Value = [double]"116993.320558452"
$ws.Range('E19').Value = [double]"56515.786358056401"
$ws.Range('F19').Value = [double]"6.2287371649486296"
$ws.Range('G19').Value = [double]"1247.91931626427"
$ws.Range('H19').Value = [double]"8.98915089460146"
$ws.Range('I19').Value = [double]"65.307000000000002"
$ws.Range('J19').Value = '..'
$ws.Range('K19').Value = '..'
$ws.Range('D20').Value = [double]"25801.1618279135"
$ws.Range('E20').Value = [double]"12269.393341150801"
$ws.Range('F20').Value = [double]"-3.1869104987660002"
$ws.Range('G20').Value = [double]"17765.112974085801"
$ws.Range('H20').Value = [double]"-9.3652863881600004E-2"
$ws.Range('I20').Value = [double]"54.085000000000001"
$ws.Range('D21').Value = [double]"35684.710453441898"
$ws.Range('E21').Value = [double]"20242.226084819999"
$ws.Range('F21').Value = [double]"1.4837407016796"
$ws.Range('G21').Value = [double]"16716.576256942899"
$ws.Range('H21').Value = [double]"1.0772871949728899"
$ws.Range('I21').Value = [double]"46.354999999999997"
$ws.Range('D22').Value = [double]"989.00266047318905"
$ws.Range('E22').Value = [double]"530.84570473784197"
$ws.Range('F22').Value = [double]"1.87878636278047"
$ws.Range('G22').Value = [double]"4434.98950884838"
$ws.Range('H22').Value = [double]"8.1377703837872293"
$ws.Range('I22').Value = [double]"53.536999999999999"
$ws.Range('D23').Value = [double]"354818.53111700597"
$ws.Range('E23').Value = [double]"164816.22425039401"
$ws.Range('F23').Value = [double]"3.2309026978242299"
$ws.Range('G23').Value = [double]"2198.4680879939701"
$ws.Range('H23').Value = [double]"2.2619561486631099"
$ws.Range('I23').Value = [double]"62.807444444444499"
$ws.Range('J23').Value = [double]"40.28"
$ws.Range('K23').Value = [double]"27.96"
$ws.Range('D24').Value = [double]"2927.0475817644301"
$ws.Range('E24').Value = [double]"1285.3028860920199"
$ws.Range('F24').Value = [double]"2.11219558419391"
$ws.Range('G24').Value = [double]"3176.9261829542902"
$ws.Range('H24').Value = [double]"-1.5189063875899999E-2"
$ws.Range('I24').Value = [double]"43.677"
$ws.Range('D25').Value = [double]"5993.7953743608696"
$ws.Range('E25').Value = [double]"3372.30318665919"
$ws.Range('G25').Value = [double]"5980.6558734069904"
$ws.Range('H25').Value = [double]"1.17976375303748"
$ws.Range('I25').Value = [double]"30.751000000000001"
$ws.Range('D26').Value = [double]"6955.35042084622"
$ws.Range('E26').Value = [double]"2211.95399725579"
$ws.Range('F26').Value = [double]"2.8949251777657499"
$ws.Range('G26').Value = [double]"1931.5052543310801"
$ws.Range('H26').Value = [double]"6.6421518407489204"
$ws.Range('I26').Value = [double]"77.046000000000006"
$ws.Range('D27').Value = [double]"314131.53149872401"
$ws.Range('E27').Value = [double]"99269.424776251704"
$ws.Range('G27').Value = [double]"3082.74319429562"
$ws.Range('H27').Value = [double]"26.7940034329601"
$ws.Range('I27').Value = [double]"79.802000000000007"
$ws.Range('D28').Value = [double]"277501.449322559"
$ws.Range('E28').Value = [double]"110347.050865439"
$ws.Range('G28').Value = [double]"5572.2489479217102"
$ws.Range('H28').Value = [double]"6.10934807127735"
$ws.Range('I28').Value = [double]"73.224000000000004"
$ws.Range('D29').Value = [double]"47298.034457627102"
$ws.Range('E29').Value = [double]"14554.7541152105"
$ws.Range('F29').Value = [double]"5.7396157374235397"
$ws.Range('G29').Value = [double]"1678.5447674649399"
$ws.Range('H29').Value = [double]"5.8194278710197196"
$ws.Range('I29').Value = [double]"84.712999999999994"
$ws.Range('J29').Value = '..'
$ws.Range('K29').Value = '..'
$ws.Range('D30').Value = [double]"29235.0600593239"
$ws.Range('E30').Value = [double]"11478.256489633501"
$ws.Range('F30').Value = [double]"3.4991246315849902"
$ws.Range('G30').Value = [double]"23131.052798447199"
$ws.Range('H30').Value = [double]"4.0388413039330802"
$ws.Range('I30').Value = [double]"58.031999999999996"
$ws.Range('D31').Value = [double]"33256.879377898003"
$ws.Range('E31').Value = [double]"11067.036252370999"
$ws.Range('G31').Value = [double]"2566.6215114851602"
$ws.Range('H31').Value = [double]"0.82795488623743996"
$ws.Range('I31').Value = [double]"56.97"
$ws.Range('D32').Value = [double]"3150.1191284496999"
$ws.Range('E32').Value = [double]"1457.1145930734799"
$ws.Range('F32').Value = [double]"7.9136979958662597"
$ws.Range('G32').Value = [double]"32128.871083076501"
$ws.Range('H32').Value = [double]"9.7704371245567891"
$ws.Range('D33').Value = [double]"18987.9606793311"
$ws.Range('G33').Value = [double]"1250.13111616022"
$ws.Range('J33').Value = '..'
$ws.Range('K33').Value = '..'
$ws.Range('D34').Value = [double]"6192.7217701484597"
$ws.Range('E34').Value = [double]"5935.0833610995896"
$ws.Range('F34').Value = [double]"5.32918198855144"
$ws.Range('G34').Value = [double]"436.76103296290302"
$ws.Range('H34').Value = [double]"30.228504863262501"
$ws.Range('I34').Value = [double]"69.891999999999996"
$ws.Range('D35').Value = [double]"195368.61746448599"
$ws.Range('E35').Value = [double]"34738.390738516697"
$ws.Range('F35').Value = [double]"0.50000000000033995"
$ws.Range('G35').Value = [double]"4293.9628663024096"
$ws.Range('H35').Value = [double]"359.09214098877902"
$ws.Range('I35').Value = [double]"48.149000000000001"
$ws.Range('D36').Value = [double]"186351.87426177901"
$ws.Range('E36').Value = [double]"69938.069501207996"
$ws.Range('F36').Value = [double]"4.94665995282186"
$ws.Range('G36').Value = [double]"3119.91916365313"
$ws.Range('H36').Value = [double]"3.6899692502561998"
$ws.Range('I36').Value = [double]"80.39"
$ws.Range('D37').Value = [double]"117880.38453202001"
$ws.Range('E37').Value = [double]"42940.920511488897"
$ws.Range('F37').Value = [double]"6.0037077530066103"
$ws.Range('G37').Value = [double]"2776.24443181795"
$ws.Range('H37').Value = [double]"2.207239312914"
$ws.Range('I37').Value = [double]"69.338999999999999"
$ws.Range('D38').Value = [double]"1245230.82592932"
$ws.Range('E38').Value = [double]"416223.66127430002"
$ws.Range('F38').Value = [double]"5.3689830287380502"
$ws.Range('G38').Value = [double]"3320.8249948869002"
$ws.Range('H38').Value = [double]"6.10934807127735"
$ws.Range('I38').Value = [double]"61.976307692307699"
$ws.Range('J38').Value = [double]"39.709090909090897"
$ws.Range('K38').Value = [double]"28.763636363636401"
$ws.Range('D39').Value = [double]"529090.85073440801"
$ws.Range('E39').Value = [double]"163137.87825924001"
$ws.Range('F39').Value = [double]"3.39999999999994"
$ws.Range('G39').Value = [double]"11869.215230997999"
$ws.Range('I39').Value = [double]"39.348999999999997"
$ws.Range('D40').Value = [double]"1467996.81495416"
$ws.Range('E40').Value = [double]"423300.33352745097"
$ws.Range('F40').Value = [double]"3.30662989992154"
$ws.Range('G40').Value = [double]"14378.0295294237"
$ws.Range('H40').Value = [double]"4.4997240400532803"
$ws.Range('I40').Value = [double]"41.54"
$ws.Range('J40').Value = [double]"31.9"
$ws.Range('K40').Value = [double]"1.5"
$ws.Range('D41').Value = [double]"147885.52744283501"
$ws.Range('E41').Value = [double]"39007.482152681099"
$ws.Range('F41').Value = [double]"28.334752000902999"
$ws.Range('G41').Value = [double]"22039.529377400799"
$ws.Range('H41').Value = [double]"2.8795052842700799"
$ws.Range('I41').Value = [double]"45.969000000000001"
$ws.Range('D42').Value = [double]"27061.934938517301"
$ws.Range('E42').Value = [double]"9891.5557777487193"
$ws.Range('F42').Value = [double]"2.4450405298867199"
$ws.Range('G42').Value = [double]"6382.4882646935503"
$ws.Range('H42').Value = [double]"3.5583029787009002"
$ws.Range('I42').Value = [double]"40.591000000000001"
$ws.Range('D43').Value = [double]"334714.84585101"
$ws.Range('E43').Value = [double]"142866.58312469599"
$ws.Range('F43').Value = [double]"7.9296679419552198"
$ws.Range('G43').Value = [double]"9217.4936207696901"
$ws.Range('H43').Value = [double]"1.39965842550414"
$ws.Range('I43').Value = [double]"43.456000000000003"
$ws.Range('D44').Value = [double]"140672.540097424"
$ws.Range('E44').Value = [double]"46688.412170077303"
$ws.Range('F44').Value = [double]"4.4053620608957198"
$ws.Range('G44').Value = [double]"11704.6533454488"
$ws.Range('H44').Value = [double]"5.7067365931682303"
$ws.Range('I44').Value = [double]"44.554000000000002"
$ws.Range('D45').Value = [double]"2647422.5140183601"
$ws.Range('E45').Value = [double]"824892.24501189496"
$ws.Range('F45').Value = [double]"5.3574332412243004"
$ws.Range('G45').Value = [double]"12882.443529284201"
$ws.Range('H45').Value = [double]"4.4997240400532803"
$ws.Range('I45').Value = [double]"42.576500000000003"
$ws.Range('J45').Value = [double]"34.200000000000003"
$ws.Range('K45').Value = [double]"2.375"
$ws.Range('D46').Value = [double]"47635.6829172465"
$ws.Range('E46').Value = [double]"17699.3822576602"
$ws.Range('F46').Value = [double]"7.1554516458795403"
$ws.Range('G46').Value = [double]"3652.85630732339"
$ws.Range('I46').Value = [double]"62.081000000000003"
$ws.Range('D47').Value = [double]"53435.002540109301"
$ws.Range('E47').Value = [double]"19748.086114474099"
$ws.Range('F47').Value = [double]"6.9045617450555401"
$ws.Range('G47').Value = [double]"2413.0113952276502"
$ws.Range('H47').Value = [double]"3.9081280003847101"
$ws.Range('I47').Value = [double]"64.756"
$ws.Range('D48').Value = [double]"4349.2278723746804"
$ws.Range('E48').Value = [double]"2112.3478342154799"
$ws.Range('F48').Value = [double]"6.9999999999995897"
$ws.Range('G48').Value = [double]"7718.9788065580397"
$ws.Range('H48').Value = [double]"1.8634419465745999"
$ws.Range('I48').Value = [double]"54.249000000000002"
$ws.Range('D49').Value = [double]"162143.42402191699"
$ws.Range('E49').Value = [double]"71712.031933781196"
$ws.Range('F49').Value = [double]"6.9999999999999902"
$ws.Range('G49').Value = [double]"5862.2976987328202"
$ws.Range('H49').Value = [double]"4.1591965791182801"
$ws.Range('I49').Value = [double]"64.064999999999998"
$ws.Range('D50').Value = [double]"6049.0569220285297"
$ws.Range('E50').Value = [double]"2034.44143273645"
$ws.Range('F50').Value = [double]"4.2654920880355496"
$ws.Range('G50').Value = [double]"2427.3636987630398"
$ws.Range('H50').Value = [double]"7.3703606115290299"
$ws.Range('I50').Value = [double]"60.271999999999998"
$ws.Range('J50').Value = [double]"38.799999999999997"
$ws.Range('K50').Value = [double]"17.2"
$ws.Range('D51').Value = [double]"196916.10929333101"
$ws.Range('E51').Value = [double]"79156.955298471803"
$ws.Range('F51').Value = [double]"5.3564779177171404"
$ws.Range('G51').Value = [double]"6272.3699092448196"
$ws.Range('H51').Value = [double]"9.9756858931185306"
$ws.Range('I51').Value = [double]"68.376999999999995"
$ws.Range('K51').Value = [double]"25.2"
$ws.Range('D52').Value = [double]"39573.454833086696"
$ws.Range('E52').Value = [double]"16155.314984516801"
$ws.Range('F52').Value = [double]"4.2917921205071403"
$ws.Range('G52').Value = [double]"2764.17899229657"
$ws.Range('H52').Value = [double]"12.596573734347301"
$ws.Range('I52').Value = [double]"52.378"
$ws.Range('D53').Value = [double]"4951.3939472652801"
$ws.Range('E53').Value = [double]"1725.46233065703"
$ws.Range('F53').Value = [double]"6.3999999999999897"
$ws.Range('G53').Value = [double]"2669.0138181796401"
$ws.Range('H53').Value = [double]"3.2739505920041498"
$ws.Range('I53').Value = [double]"55.424999999999997"
$ws.Range('D54').Value = [double]"7993.5650538333402"
$ws.Range('E54').Value = [double]"3508.8634526778201"
$ws.Range('F54').Value = [double]"5.0109751952555301"
$ws.Range('G54').Value = [double]"1543.1592768017999"
$ws.Range('H54').Value = [double]"7.8158724211053503"
$ws.Range('I54').Value = [double]"75.951999999999998"
$ws.Range('D55').Value = [double]"51257.746486191099"
$ws.Range('E55').Value = [double]"19656.025198198098"
$ws.Range('F55').Value = [double]"3.0548521639981501"
$ws.Range('G55').Value = [double]"2338.9398883797799"
$ws.Range('H55').Value = [double]"3.80833257901772"
$ws.Range('I55').Value = [double]"65.92"
$ws.Range('D56').Value = [double]"33067.849763874903"
$ws.Range('E56').Value = [double]"14922.728893445201"
$ws.Range('F56').Value = [double]"1.3993039513252401"
$ws.Range('G56').Value = [double]"1315.82868013702"
$ws.Range('I56').Value = [double]"72.778000000000006"
$ws.Range('D57').Value = [double]"1159167.65372629"
$ws.Range('E57').Value = [double]"441424.28367733402"
$ws.Range('G57').Value = [double]"5483.2647609343903"
$ws.Range('H57').Value = [double]"16.952790783419601"
$ws.Range('I57').Value = [double]"58.311"
$ws.Range('D58').Value = [double]"65095.509147623103"
$ws.Range('E58').Value = [double]"27639.683564947201"
$ws.Range('F58').Value = [double]"6.0667708633988902"
$ws.Range('G58').Value = [double]"3783.7888412477801"
$ws.Range('H58').Value = [double]"2.1784077383737399"
$ws.Range('I58').Value = [double]"50.228000000000002"
$ws.Range('D59').Value = [double]"15003.347045422999"
$ws.Range('E59').Value = [double]"4147.8705956497397"
$ws.Range('G59').Value = [double]"1842.8589786013199"
$ws.Range('H59').Value = [double]"11.8739588840234"
$ws.Range('I59').Value = [double]"53.298000000000002"
$ws.Range('D60').Value = [double]"20268.4232884501"
$ws.Range('E60').Value = [double]"8430.3896229531001"
$ws.Range('F60').Value = [double]"5.2588999999999997"
$ws.Range('G60').Value = [double]"2341.1327447572598"
$ws.Range('H60').Value = [double]"4.5483174278751797"
$ws.Range('I60').Value = [double]"57.304000000000002"
$ws.Range('D61').Value = [double]"1866907.4468590501"
$ws.Range('E61').Value = [double]"730073.86719171796"
$ws.Range('F61').Value = [double]"4.3875039817292798"
$ws.Range('G61').Value = [double]"4551.3434618502497"
$ws.Range('H61').Value = [double]"16.952790783419601"
$ws.Range('I61').Value = [double]"61.0262666666667"
$ws.Range('J61').Value = [double]"38.093333333333298"
$ws.Range('K61').Value = [double]"22.113333333333301"
$ws.Range('D62').Value = [double]"7479684.8719005296"
$ws.Range('E62').Value = [double]"2754427.4616089002"
$ws.Range('F62').Value = [double]"4.8523047500569403"
$ws.Range('G62').Value = [double]"5541.7256151019601"
$ws.Range('H62').Value = [double]"4.5556051808843199"
$ws.Range('I62').Value = [double]"60.013905660377397"
$ws.Range('J62').Value = [double]"41.6111111111111"
$ws.Range('K62').Value = [double]"26.32"
$ws.Range('D63').Value = [double]"138168752.01852101"
$ws.Range('E63').Value = [double]"92227592.093919605"
$ws.Range('F63').Value = [double]"6.3546695842029601"
$ws.Range('G63').Value = [double]"21850.4879059368"
$ws.Range('H63').Value = [double]"3.2123454846075199"
$ws.Range('I63').Value = [double]"59.403128000000002"
$ws.Range('J63').Value = [double]"34.991578947368403"
$ws.Range('K63').Value = [double]"2.34"
$ws.Range('D64').Value = [double]"10763828.400106501"
$ws.Range('E64').Value = [double]"5049232.7616223805"
$ws.Range('F64').Value = [double]"7.0118886659601198"
$ws.Range('G64').Value = [double]"16782.8874556216"
$ws.Range('H64').Value = [double]"5.6929663312506698"
$ws.Range('I64').Value = [double]"60.535103448275898"
$ws.Range('J64').Value = [double]"45.461111111111101"
$ws.Range('K64').Value = [double]"3.5833333333333299"
$ws.Range('D65').Value = [double]"51850404.541612297"
$ws.Range('E65').Value = [double]"25476740.087795202"
$ws.Range('F65').Value = [double]"7.3388234896270497"
$ws.Range('G65').Value = [double]"12488.87507127"
$ws.Range('H65').Value = [double]"0.85250710422560005"
$ws.Range('I65').Value = [double]"54.919566666666697"
$ws.Range('J65').Value = [double]"35.174999999999997"
$ws.Range('K65').Value = [double]"5"
$ws.Range('D66').Value = [double]"145648436.890421"
$ws.Range('E66').Value = [double]"94982019.555528596"
$ws.Range('F66').Value = [double]"6.2775165727326696"
$ws.Range('G66').Value = [double]"18924.589148430001"
$ws.Range('H66').Value = [double]"3.2224035973881899"
$ws.Range('I66').Value = [double]"59.584988764044901"
$ws.Range('J66').Value = [double]"37.119285714285702"
$ws.Range('K66').Value = [double]"10.047857142857101"
$ws.Range('D67').Value = [double]"3086870.2936187699"
$ws.Range('E67').Value = [double]"984535.33921808004"
$ws.Range('F67').Value = [double]"5.5152235531443896"
$ws.Range('G67').Value = [double]"5262.8493452994198"
$ws.Range('H67').Value = [double]"4.4997240400532803"
$ws.Range('I67').Value = [double]"58.797350000000002"
$ws.Range('J67').Value = [double]"41"
$ws.Range('K67').Value = [double]"29.893750000000001"
$ws.Range('D68').Value = [double]"4235336.6999360798"
$ws.Range('E68').Value = [double]"1449827.78453374"
$ws.Range('F68').Value = [double]"4.8847830738179203"
$ws.Range('G68').Value = [double]"6535.5576710919704"
$ws.Range('H68').Value = [double]"4.4997240400532803"
$ws.Range('I68').Value = [double]"54.550800000000002"
$ws.Range('J68').Value = [double]"37.576190476190497"
$ws.Range('K68').Value = [double]"18.709523809523802"
$ws.Range('D69').Value = [double]"631078.39001443796"
$ws.Range('E69').Value = [double]"243578.970419088"
$ws.Range('F69').Value = [double]"6.5619510723797196"
$ws.Range('G69').Value = [double]"3365.5701821692301"
$ws.Range('H69').Value = [double]"3.6899692502561998"
$ws.Range('I69').Value = [double]"71.383666666666699"
$ws.Range('D70').Value = [double]"612731.02252825198"
$ws.Range('E70').Value = [double]"250673.82452555199"
$ws.Range('F70').Value = [double]"2.8614956668385201"
$ws.Range('G70').Value = [double]"2909.3772923958099"
$ws.Range('H70').Value = [double]"8.98915089460146"
$ws.Range('I70').Value = [double]"63.5178181818182"
$ws.Range('J70').Value = [double]"42.342857142857198"
$ws.Range('K70').Value = [double]"31.842857142857099"
$ws.Range('D71').Value = [double]"1866907.4468590501"
$ws.Range('E71').Value = [double]"730073.86719171796"
$ws.Range('F71').Value = [double]"4.3875039817292798"
$ws.Range('G71').Value = [double]"4551.3434618502497"
$ws.Range('H71').Value = [double]"16.952790783419601"
$ws.Range('I71').Value = [double]"61.0262666666667"
$ws.Range('J71').Value = [double]"38.093333333333298"
$ws.Range('K71').Value = [double]"22.113333333333301"
$ws.Range('D72').Value = [double]"943011.81106247497"
$ws.Range('E72').Value = [double]"306443.12743671099"
$ws.Range('F72').Value = [double]"5.2993014559947698"
$ws.Range('G72').Value = [double]"3461.18688281184"
$ws.Range('H72').Value = [double]"26.7940034329601"
$ws.Range('I72').Value = [double]"60.238750000000003"
$ws.Range('J72').Value = [double]"39.733333333333299"
$ws.Range('K72').Value = [double]"33.383333333333297"
$ws.Range('D73').Value = [double]"1751261.0100242"
$ws.Range('E73').Value = [double]"773650.74782387295"
$ws.Range('F73').Value = [double]"4.6496045285146703"
$ws.Range('G73').Value = [double]"4616.0132879095199"
$ws.Range('H73').Value = [double]"4.5556051808843199"
$ws.Range('I73').Value = [double]"64.736266666666694"
$ws.Range('J73').Value = [double]"48.628571428571398"
$ws.Range('K73').Value = [double]"32.221428571428604"
$ws.Range('D74').Value = [double]"1179425.69906419"
$ws.Range('E74').Value = [double]"401591.91148444399"
$ws.Range('F74').Value = [double]"7.9100084271970204"
$ws.Range('G74').Value = [double]"11324.1857831082"
$ws.Range('H74').Value = [double]"5.7067365931682303"
$ws.Range('I74').Value = [double]"42.783799999999999"
$ws.Range('D75').Value = [double]"278367.07256546698"
$ws.Range('E75').Value = [double]"94935.970828012301"
$ws.Range('F75').Value = [double]"1.4742448067929499"
$ws.Range('G75').Value = [double]"4018.6429615070601"
$ws.Range('H75').Value = [double]"25.765270173944401"
$ws.Range('I75').Value = [double]"63.5518"
$ws.Range('D76').Value = [double]"9127759.4906712696"
$ws.Range('E76').Value = [double]"3378212.6471554399"
$ws.Range('F76').Value = [double]"3.1673768086996299"
$ws.Range('G76').Value = [double]"13779.985251858399"
$ws.Range('H76').Value = [double]"1.5600797515304301"
$ws.Range('I76').Value = [double]"65.623199999999997"
$ws.Range('J76').Value = [double]"37.314285714285703"
$ws.Range('K76').Value = [double]"2.3285714285714301"
$ws.Range('D77').Value = [double]"7122809.0723406496"
$ws.Range('E77').Value = [double]"3310448.0174690099"
$ws.Range('F77').Value = [double]"7.4975057224951103"
$ws.Range('G77').Value = [double]"16437.421429324299"
$ws.Range('H77').Value = [double]"8.3015736896473697"
$ws.Range('I77').Value = [double]"59.862166666666702"
$ws.Range('J77').Value = [double]"44.655555555555601"
$ws.Range('K77').Value = [double]"2.6"
$ws.Range('D78').Value = [double]"21897662.890040401"
$ws.Range('E78').Value = [double]"17193799.329353102"
$ws.Range('F78').Value = [double]"5.5688233784039696"
$ws.Range('G78').Value = [double]"49214.207657723397"
$ws.Range('H78').Value = [double]"2.84249072759921"
$ws.Range('I78').Value = [double]"58.454925925925899"
$ws.Range('J78').Value = [double]"30.766666666666701"
$ws.Range('K78').Value = [double]"0.31111111111111001"
$ws.Range('D79').Value = [double]"68124121.305473804"
$ws.Range('E79').Value = [double]"58455215.938100196"
$ws.Range('F79').Value = [double]"5.7436413895197296"
$ws.Range('G79').Value = [double]"49738.131253412597"
$ws.Range('H79').Value = [double]"3.4012777608766598"
$ws.Range('I79').Value = [double]"60.308236842105302"
$ws.Range('J79').Value = [double]"33.089189189189199"
$ws.Range('K79').Value = [double]"0.57297297297297001"
$ws.Range('D80').Value = [double]"2202725.2305942001"
$ws.Range('E80').Value = [double]"786806.77821717504"
$ws.Range('F80').Value = [double]"4.7815438225206899"
$ws.Range('G80').Value = [double]"6718.5897365742103"
$ws.Range('H80').Value = [double]"16.952790783419601"
$ws.Range('I80').Value = [double]"54.186399999999999"
$ws.Range('D81').Value = [double]"11489164.9881438"
$ws.Range('E81').Value = [double]"4713062.6907911403"
$ws.Range('F81').Value = [double]"4.7158118531321502"
$ws.Range('G81').Value = [double]"23211.351773953698"
$ws.Range('H81').Value = [double]"6.6944554174957798"
$ws.Range('I81').Value = [double]"59.679136363636403"
$ws.Range('J81').Value = [double]"34.325000000000003"
$ws.Range('K81').Value = [double]"6.1875"
$ws.Range('D82').Value = [double]"5276959.6413063304"
$ws.Range('E82').Value = [double]"1967620.6833917301"
$ws.Range('F82').Value = [double]"4.8818420008640597"
$ws.Range('G82').Value = [double]"5174.5002579084103"
$ws.Range('H82').Value = [double]"4.4997240400532803"
$ws.Range('I82').Value = [double]"61.3691395348838"
$ws.Range('J82').Value = [double]"41.938461538461503"
$ws.Range('K82').Value = [double]"27.638461538461499"
$ws.Range('D83').Value = [double]"126679587.030377"
$ws.Range('E83').Value = [double]"87514529.403128505"
$ws.Range('F83').Value = [double]"6.5033052666026299"
$ws.Range('G83').Value = [double]"21734.264565233902"
$ws.Range('H83').Value = [double]"3.2123454846075199"
$ws.Range('I83').Value = [double]"59.3441747572816"
$ws.Range('J83').Value = [double]"35.052873563218398"
$ws.Range('K83').Value = [double]"1.9862068965517199"
$ws.Range('D84').Value = [double]"1207530.8992618001"
$ws.Range('E84').Value = [double]"411166.26219836099"
$ws.Range('F84').Value = [double]"4.5636666343181602"
$ws.Range('G84').Value = [double]"2210.2992181843201"
$ws.Range('H84').Value = [double]"8.98915089460146"
$ws.Range('I84').Value = [double]"65.086565217391296"
$ws.Range('J84').Value = [double]"39.200000000000003"
$ws.Range('K84').Value = [double]"36.938888888888897"
$ws.Range('D85').Value = [double]"62015.771661434497"
$ws.Range('E85').Value = [double]"16956.521156046001"
$ws.Range('F85').Value = [double]"-0.99999999999989997"
$ws.Range('G85').Value = [double]"1899.93700741524"
$ws.Range('H85').Value = [double]"26.0218208123664"
$ws.Range('I85').Value = [double]"49.044499999999999"
$ws.Range('D86').Value = [double]"5089904.5974087697"
$ws.Range('E86').Value = [double]"1808824.31708762"
$ws.Range('F86').Value = [double]"4.2509028253117904"
$ws.Range('G86').Value = [double]"7132.63222266304"
$ws.Range('H86').Value = [double]"5.7067365931682303"
$ws.Range('I86').Value = [double]"56.684869565217397"
$ws.Range('J86').Value = [double]"42.123809523809499"
$ws.Range('K86').Value = [double]"22.1380952380952"
$ws.Range('D87').Value = [double]"22307340.2729045"
$ws.Range('E87').Value = [double]"6865905.0118653905"
$ws.Range('F87').Value = [double]"6.43821276707693"
$ws.Range('G87').Value = [double]"8428.4202352984703"
$ws.Range('H87').Value = [double]"5.5069937221637497"
$ws.Range('I87').Value = [double]"58.891517241379297"
$ws.Range('J87').Value = [double]"35.908333333333303"
$ws.Range('K87').Value = [double]"5.4666666666666703"
$ws.Range('D88').Value = [double]"1179099.25610151"
$ws.Range('E88').Value = [double]"532979.76772985503"
$ws.Range('F88').Value = [double]"7.7358398673024302"
$ws.Range('G88').Value = [double]"15368.023591216999"
$ws.Range('H88').Value = [double]"4.5556051808843199"
$ws.Range('I88').Value = [double]"54.284857142857199"
$ws.Range('D89').Value = [double]"49768341.708983101"
$ws.Range('E89').Value = [double]"26917735.884856202"
$ws.Range('F89').Value = [double]"7.6245635184190297"
$ws.Range('G89').Value = [double]"20422.2509312533"
$ws.Range('H89').Value = [double]"0.85250710422560005"
$ws.Range('I89').Value = [double]"57.515142857142898"
$ws.Range('J89').Value = [double]"38.421875"
$ws.Range('K89').Value = [double]"1.8968750000000001"
$ws.Range('D90').Value = [double]"65865048.972489901"
$ws.Range('E90').Value = [double]"58371303.927981898"
$ws.Range('F90').Value = [double]"5.3889823405144703"
$ws.Range('G90').Value = [double]"56048.667841997398"
$ws.Range('H90').Value = [double]"3.2123454846075199"
$ws.Range('I90').Value = [double]"62.335306122448998"
$ws.Range('J90').Value = [double]"31.494871794871798"
$ws.Range('K90').Value = [double]"0.28461538461538"
$ws.Range('D91').Value = [double]"1843786.78031773"
$ws.Range('E91').Value = [double]"641011.35125858197"
$ws.Range('F91').Value = [double]"4.2570995351535199"
$ws.Range('G91').Value = [double]"2624.3349801566201"
$ws.Range('H91').Value = [double]"8.98915089460146"
$ws.Range('I91').Value = [double]"62.415727272727302"
$ws.Range('J91').Value = [double]"40.553571428571402"
$ws.Range('K91').Value = [double]"32.9892857142857"
$ws.Range('D92').Value = [double]"1804536.1906908399"
$ws.Range('E92').Value = [double]"608717.24431728804"
$ws.Range('F92').Value = [double]"2.5839646968124201"
$ws.Range('G92').Value = [double]"5635.6652049346703"
$ws.Range('H92').Value = [double]"5.5576564158240904"
$ws.Range('I92').Value = [double]"58.820181818181801"
$ws.Range('J92').Value = [double]"33.214285714285701"
$ws.Range('K92').Value = [double]"9.9142857142857199"
$ws.Range('D93').Value = [double]"45601.851249651198"
$ws.Range('E93').Value = [double]"18589.329838409401"
$ws.Range('F93').Value = [double]"4.3287792675260501"
$ws.Range('G93').Value = [double]"9238.3228848237504"
$ws.Range('H93').Value = [double]"4.0388413039330802"
$ws.Range('I93').Value = [double]"52.984000000000002"
$ws.Range('D94').Value = [double]"1192957.35004167"
$ws.Range('E94').Value = [double]"698590.62040724303"
$ws.Range('F94').Value = [double]"8.2832601386794593"
$ws.Range('G94').Value = [double]"24543.115213342098"
$ws.Range('H94').Value = [double]"2.3048595904049201"
$ws.Range('I94').Value = [double]"61.475227272727302"
$ws.Range('J94').Value = [double]"35.118181818181803"
$ws.Range('K94').Value = [double]"5.7545454545454602"
$ws.Range('D95').Value = [double]"850627.53228690696"
$ws.Range('E95').Value = [double]"327914.35472718102"
$ws.Range('F95').Value = [double]"5.9190857560761696"
$ws.Range('G95').Value = [double]"2575.7272846003302"
$ws.Range('H95').Value = [double]"22.0207676245778"
$ws.Range('I95').Value = [double]"66.161937499999993"
$ws.Range('D96').Value = [double]"1766049.81869214"
$ws.Range('E96').Value = [double]"597964.38238594402"
$ws.Range('F96').Value = [double]"5.1624692004160497"
$ws.Range('G96').Value = [double]"11383.6308582742"
$ws.Range('H96').Value = [double]"8.0018662835986305"
$ws.Range('I96').Value = [double]"56.045124999999999"
$ws.Range('K96').Value = [double]"2"
$ws.Range('D97').Value = [double]"3617062.9924139101"
$ws.Range('E97').Value = [double]"1344604.63925397"
$ws.Range('F97').Value = [double]"5.2506340834632201"
$ws.Range('G97').Value = [double]"3654.9202197483201"
$ws.Range('H97').Value = [double]"16.952790783419601"
$ws.Range('I97').Value = [double]"62.6536756756757"
$ws.Range('J97').Value = [double]"41.424137931034501"
$ws.Range('K97').Value = [double]"32.0724137931035"
$ws.Range('D98').Value = [double]"5425606.8276800402"
$ws.Range('E98').Value = [double]"1632183.0259404799"
$ws.Range('F98').Value = [double]"4.64333534714069"
$ws.Range('G98').Value = [double]"7569.57103499401"
$ws.Range('H98').Value = [double]"8.9014616352852904"
$ws.Range('I98').Value = [double]"55.262631578947399"
$ws.Range('J98').Value = [double]"38.58"
$ws.Range('K98').Value = [double]"8.0500000000000007"

# --- A104 previously held the stale "Source:" footnote; upstream removed that footnote, leaving the cell as #N/A ---
$ws.Range('A104').Value = '#N/A'
